# Eliminate colonne come da indicazioni
# S1#111EUROMEDICAL / EUROMEDICAL / MEDIWEB / 1.0
#
# This clears the "ERRORE BLOCCANTE -> Tipo Documento non gestito" marker
# columns (J and K) that were previously filled in with "NO" /
# "Tipo Documento non gestito" on every test-case row, and resets the
# sheet's frozen-pane / active-cell selection back to the top of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,21,23,24,25,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,105,106,107,108,109)

foreach ($r in $rows) {
    $ws.Range("J$r").Value = ""
    $ws.Range("K$r").Value = ""
}

# The very first cleared cell keeps a stray single-space value, matching
# the author's actual edit (everything else ends up fully blank).
$ws.Range("J8").Value = " "

# Reset the view: frozen pane top-left cell and the active selection.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("E5").Select()
